$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2852.45616702906
$ws.Range("B3").Value = 3222.810395715689
$ws.Range("B4").Value = 3760.41217101072
$ws.Range("B5").Value = 4215.576617825329
$ws.Range("B6").Value = 4552.606182681607
$ws.Range("B7").Value = 4837.553047886875
$ws.Range("B8").Value = 5105.324759595394
$ws.Range("B9").Value = 5327.087680414962
$ws.Range("B10").Value = 5543.973262663902
$ws.Range("B11").Value = 5722.967529782031
$ws.Range("B12").Value = 5791.712971149964
$ws.Range("B13").Value = 6009.995218597167
$ws.Range("B14").Value = 6138.979774853446
$ws.Range("B15").Value = 6221.394152096422
$ws.Range("B16").Value = 6351.835669819631
$ws.Range("B17").Value = 6389.913579031808
$ws.Range("B18").Value = 6457.899256995265
$ws.Range("B19").Value = 6493.189068226568
$ws.Range("B20").Value = 6573.623914589975
$ws.Range("B21").Value = 6595.358086449128
$ws.Range("B22").Value = 6622.658099770334
$ws.Range("B23").Value = 6616.094729561163
$ws.Range("B24").Value = 6631.073662321366
$ws.Range("B25").Value = 6631.255821965488
$ws.Range("B26").Value = 6594.344551836674
$ws.Range("B27").Value = 6571.935987152133
$ws.Range("B28").Value = 6513.163518606098
$ws.Range("B29").Value = 6489.297511363043
$ws.Range("B30").Value = 6405.067096257225
$ws.Range("B31").Value = 6360.456442666165
$ws.Range("B32").Value = 6340.285597746446
$ws.Range("B33").Value = 6230.041707158317
$ws.Range("B34").Value = 6111.311892980305
$ws.Range("B35").Value = 6046.065743282183
$ws.Range("B36").Value = 5939.768889263428
$ws.Range("B37").Value = 5874.682892777044
$ws.Range("B38").Value = 5706.213541061077
$ws.Range("B39").Value = 5572.991581891588
$ws.Range("B40").Value = 5427.926963202865
$ws.Range("B41").Value = 5275.17731616132
$ws.Range("B42").Value = 5159.035857758249
$ws.Range("B43").Value = 5031.489641223315
$ws.Range("B44").Value = 4891.717422088948
$ws.Range("B45").Value = 4756.95698502396
$ws.Range("B46").Value = 4684.089424761425
$ws.Range("B47").Value = 4534.829679974262
$ws.Range("B48").Value = 4359.579989905586
$ws.Range("B49").Value = 4192.955161297428
$ws.Range("B50").Value = 4044.558644583612
$ws.Range("B51").Value = 3864.944845544914
$ws.Range("B52").Value = 3668.648383986862
$ws.Range("B53").Value = 3500.895877670191
$ws.Range("B54").Value = 3355.983223828149
$ws.Range("B55").Value = 3145.085644519123
$ws.Range("B56").Value = 2890.445413936145
$ws.Range("B57").Value = 2741.557846218197
$ws.Range("B58").Value = 2603.40694896459
$ws.Range("B59").Value = 2470.387404265228
$ws.Range("B60").Value = 2402.454256589491
$ws.Range("B61").Value = 2351.931898994147
$ws.Range("B62").Value = 2320.473191916665
